# Auto-generated edit script applying the Gungnir_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1212.4231  # ALC!H38
$ws.Cells.Item(38, 9).Value = 198.11111  # ALC!I38
$ws.Cells.Item(38, 10).Value = 1749.4117  # ALC!J38
$ws.Cells.Item(38, 11).Value = 594.3333299999999  # ALC!K38
$ws.Cells.Item(38, 12).Value = 5248.2351  # ALC!L38
$ws.Cells.Item(38, 13).Value = -222.3333299999999  # ALC!M38
$ws.Cells.Item(38, 14).Value = -5992.2351  # ALC!N38

$ws.Cells.Item(62, 8).Value = 22740256  # ALC!H62
$ws.Cells.Item(62, 9).Value = 25013780  # ALC!I62
$ws.Cells.Item(62, 10).Value = 5000  # ALC!J62
$ws.Cells.Item(62, 11).Value = 25013780  # ALC!K62
$ws.Cells.Item(62, 12).Value = 5000  # ALC!L62
$ws.Cells.Item(62, 13).Value = -25013156  # ALC!M62
$ws.Cells.Item(62, 14).Value = -6248  # ALC!N62

$ws.Cells.Item(65, 8).Value = 22740256  # ALC!H65
$ws.Cells.Item(65, 9).Value = 25013780  # ALC!I65
$ws.Cells.Item(65, 10).Value = 5000  # ALC!J65
$ws.Cells.Item(65, 11).Value = 125068900  # ALC!K65
$ws.Cells.Item(65, 12).Value = 25000  # ALC!L65
$ws.Cells.Item(65, 13).Value = -125065780  # ALC!M65
$ws.Cells.Item(65, 14).Value = -31240  # ALC!N65

$ws.Cells.Item(86, 8).Value = 5815.75  # ALC!H86
$ws.Cells.Item(86, 9).Value = 7790.2144  # ALC!I86
$ws.Cells.Item(86, 10).Value = 1208.6666  # ALC!J86
$ws.Cells.Item(86, 11).Value = 7790.2144  # ALC!K86
$ws.Cells.Item(86, 12).Value = 1208.6666  # ALC!L86
$ws.Cells.Item(86, 13).Value = -6667.2144  # ALC!M86
$ws.Cells.Item(86, 14).Value = -3454.6666  # ALC!N86

$ws.Cells.Item(89, 8).Value = 5815.75  # ALC!H89
$ws.Cells.Item(89, 9).Value = 7790.2144  # ALC!I89
$ws.Cells.Item(89, 10).Value = 1208.6666  # ALC!J89
$ws.Cells.Item(89, 11).Value = 38951.072  # ALC!K89
$ws.Cells.Item(89, 12).Value = 6043.333000000001  # ALC!L89
$ws.Cells.Item(89, 13).Value = -33335.072  # ALC!M89
$ws.Cells.Item(89, 14).Value = -17275.333  # ALC!N89

$ws.Cells.Item(113, 8).Value = 3415.158  # ALC!H113
$ws.Cells.Item(113, 9).Value = 3412.5715  # ALC!I113
$ws.Cells.Item(113, 10).Value = 3416.6667  # ALC!J113
$ws.Cells.Item(113, 11).Value = 3412.5715  # ALC!K113
$ws.Cells.Item(113, 12).Value = 3416.6667  # ALC!L113
$ws.Cells.Item(113, 13).Value = -158.5715  # ALC!M113
$ws.Cells.Item(113, 14).Value = -9924.6667  # ALC!N113

$ws.Cells.Item(132, 8).Value = 12507472  # ALC!H132
$ws.Cells.Item(132, 9).Value = 15632110  # ALC!I132
$ws.Cells.Item(132, 10).Value = 8918  # ALC!J132
$ws.Cells.Item(132, 11).Value = 46896330  # ALC!K132
$ws.Cells.Item(132, 12).Value = 26754  # ALC!L132
$ws.Cells.Item(132, 13).Value = -46893800  # ALC!M132
$ws.Cells.Item(132, 14).Value = -31814  # ALC!N132

$ws.Cells.Item(135, 8).Value = 8929558  # ALC!H135
$ws.Cells.Item(135, 9).Value = 8929558  # ALC!I135
$ws.Cells.Item(135, 10).Value = 0  # ALC!J135
$ws.Cells.Item(135, 11).Value = 80366022  # ALC!K135
$ws.Cells.Item(135, 12).Value = 0  # ALC!L135
$ws.Cells.Item(135, 13).Value = -80363487  # ALC!M135
$ws.Cells.Item(135, 14).Value = $null  # ALC!N135 (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 1990  # ARM!H25
$ws.Cells.Item(25, 9).Value = 1990  # ARM!I25
$ws.Cells.Item(25, 10).Value = 0  # ARM!J25
$ws.Cells.Item(25, 11).Value = 1990  # ARM!K25
$ws.Cells.Item(25, 12).Value = 0  # ARM!L25
$ws.Cells.Item(25, 13).Value = -1588  # ARM!M25

$ws.Cells.Item(31, 8).Value = 3000  # ARM!H31
$ws.Cells.Item(31, 9).Value = 3000  # ARM!I31
$ws.Cells.Item(31, 10).Value = 0  # ARM!J31
$ws.Cells.Item(31, 11).Value = 3000  # ARM!K31
$ws.Cells.Item(31, 12).Value = 0  # ARM!L31
$ws.Cells.Item(31, 13).Value = -2706  # ARM!M31

$ws.Cells.Item(32, 8).Value = 23535454  # ARM!H32
$ws.Cells.Item(32, 9).Value = 12663317  # ARM!I32
$ws.Cells.Item(32, 10).Value = 166685250  # ARM!J32
$ws.Cells.Item(32, 11).Value = 12663317  # ARM!K32
$ws.Cells.Item(32, 12).Value = 166685250  # ARM!L32
$ws.Cells.Item(32, 13).Value = -12663030  # ARM!M32
$ws.Cells.Item(32, 14).Value = -166685824  # ARM!N32

$ws.Cells.Item(35, 8).Value = 1000  # ARM!H35
$ws.Cells.Item(35, 9).Value = 1000  # ARM!I35
$ws.Cells.Item(35, 10).Value = 0  # ARM!J35
$ws.Cells.Item(35, 11).Value = 1000  # ARM!K35
$ws.Cells.Item(35, 12).Value = 0  # ARM!L35
$ws.Cells.Item(35, 13).Value = -594  # ARM!M35

$ws.Cells.Item(69, 8).Value = 62972.332  # ARM!H69
$ws.Cells.Item(69, 9).Value = 0  # ARM!I69
$ws.Cells.Item(69, 10).Value = 62972.332  # ARM!J69
$ws.Cells.Item(69, 11).Value = 0  # ARM!K69
$ws.Cells.Item(69, 12).Value = 62972.332  # ARM!L69
$ws.Cells.Item(69, 14).Value = -64470.332  # ARM!N69

$ws.Cells.Item(72, 8).Value = 62972.332  # ARM!H72
$ws.Cells.Item(72, 9).Value = 0  # ARM!I72
$ws.Cells.Item(72, 10).Value = 62972.332  # ARM!J72
$ws.Cells.Item(72, 11).Value = 0  # ARM!K72
$ws.Cells.Item(72, 12).Value = 188916.996  # ARM!L72
$ws.Cells.Item(72, 14).Value = -196404.996  # ARM!N72

$ws.Cells.Item(93, 8).Value = 0  # ARM!H93
$ws.Cells.Item(93, 9).Value = 0  # ARM!I93
$ws.Cells.Item(93, 10).Value = 0  # ARM!J93
$ws.Cells.Item(93, 11).Value = 0  # ARM!K93
$ws.Cells.Item(93, 12).Value = 0  # ARM!L93
$ws.Cells.Item(93, 14).Value = $null  # ARM!N93 (removed)

$ws.Cells.Item(132, 8).Value = 30307542  # ARM!H132
$ws.Cells.Item(132, 9).Value = 41667630  # ARM!I132
$ws.Cells.Item(132, 10).Value = 13966.889  # ARM!J132
$ws.Cells.Item(132, 11).Value = 125002890  # ARM!K132
$ws.Cells.Item(132, 12).Value = 41900.667  # ARM!L132
$ws.Cells.Item(132, 13).Value = -125000360  # ARM!M132
$ws.Cells.Item(132, 14).Value = -46960.667  # ARM!N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 0  # BSM!H37
$ws.Cells.Item(37, 9).Value = 0  # BSM!I37
$ws.Cells.Item(37, 10).Value = 0  # BSM!J37
$ws.Cells.Item(37, 11).Value = 0  # BSM!K37
$ws.Cells.Item(37, 12).Value = 0  # BSM!L37
$ws.Cells.Item(37, 13).Value = $null  # BSM!M37 (removed)

$ws.Cells.Item(102, 8).Value = 18015  # BSM!H102
$ws.Cells.Item(102, 9).Value = 12360.429  # BSM!I102
$ws.Cells.Item(102, 10).Value = 37806  # BSM!J102
$ws.Cells.Item(102, 11).Value = 12360.429  # BSM!K102
$ws.Cells.Item(102, 12).Value = 37806  # BSM!L102
$ws.Cells.Item(102, 13).Value = -9115.429  # BSM!M102
$ws.Cells.Item(102, 14).Value = -44296  # BSM!N102

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2725.4119  # CRP!H16
$ws.Cells.Item(16, 9).Value = 1907.4286  # CRP!I16
$ws.Cells.Item(16, 10).Value = 3298  # CRP!J16
$ws.Cells.Item(16, 11).Value = 1907.4286  # CRP!K16
$ws.Cells.Item(16, 12).Value = 3298  # CRP!L16
$ws.Cells.Item(16, 13).Value = -1620.4286  # CRP!M16
$ws.Cells.Item(16, 14).Value = -3872  # CRP!N16

$ws.Cells.Item(31, 8).Value = 974.62  # CRP!H31
$ws.Cells.Item(31, 9).Value = 787.63043  # CRP!I31
$ws.Cells.Item(31, 10).Value = 3125  # CRP!J31
$ws.Cells.Item(31, 11).Value = 787.63043  # CRP!K31
$ws.Cells.Item(31, 12).Value = 3125  # CRP!L31
$ws.Cells.Item(31, 13).Value = -492.63043  # CRP!M31
$ws.Cells.Item(31, 14).Value = -3715  # CRP!N31

$ws.Cells.Item(34, 8).Value = 974.62  # CRP!H34
$ws.Cells.Item(34, 9).Value = 787.63043  # CRP!I34
$ws.Cells.Item(34, 10).Value = 3125  # CRP!J34
$ws.Cells.Item(34, 11).Value = 787.63043  # CRP!K34
$ws.Cells.Item(34, 12).Value = 3125  # CRP!L34
$ws.Cells.Item(34, 13).Value = -585.63043  # CRP!M34
$ws.Cells.Item(34, 14).Value = -3529  # CRP!N34

$ws.Cells.Item(113, 8).Value = 2725.4119  # CRP!H113
$ws.Cells.Item(113, 9).Value = 1907.4286  # CRP!I113
$ws.Cells.Item(113, 10).Value = 3298  # CRP!J113
$ws.Cells.Item(113, 11).Value = 1907.4286  # CRP!K113
$ws.Cells.Item(113, 12).Value = 3298  # CRP!L113
$ws.Cells.Item(113, 13).Value = 262.5714  # CRP!M113
$ws.Cells.Item(113, 14).Value = -7638  # CRP!N113

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 50535.55  # CUL!H12
$ws.Cells.Item(12, 9).Value = 6.5  # CUL!I12
$ws.Cells.Item(12, 10).Value = 63167.812  # CUL!J12
$ws.Cells.Item(12, 11).Value = 19.5  # CUL!K12
$ws.Cells.Item(12, 12).Value = 189503.436  # CUL!L12
$ws.Cells.Item(12, 13).Value = 153.5  # CUL!M12
$ws.Cells.Item(12, 14).Value = -189849.436  # CUL!N12

$ws.Cells.Item(86, 8).Value = 150  # CUL!H86
$ws.Cells.Item(86, 9).Value = 120  # CUL!I86
$ws.Cells.Item(86, 10).Value = 300  # CUL!J86
$ws.Cells.Item(86, 11).Value = 360  # CUL!K86
$ws.Cells.Item(86, 12).Value = 900  # CUL!L86
$ws.Cells.Item(86, 13).Value = 826  # CUL!M86
$ws.Cells.Item(86, 14).Value = -3272  # CUL!N86

$ws.Cells.Item(89, 8).Value = 150  # CUL!H89
$ws.Cells.Item(89, 9).Value = 120  # CUL!I89
$ws.Cells.Item(89, 10).Value = 300  # CUL!J89
$ws.Cells.Item(89, 11).Value = 1080  # CUL!K89
$ws.Cells.Item(89, 12).Value = 2700  # CUL!L89
$ws.Cells.Item(89, 13).Value = 4848  # CUL!M89
$ws.Cells.Item(89, 14).Value = -14556  # CUL!N89

$ws.Cells.Item(97, 8).Value = 919  # CUL!H97
$ws.Cells.Item(97, 9).Value = 800  # CUL!I97
$ws.Cells.Item(97, 10).Value = 978.5  # CUL!J97
$ws.Cells.Item(97, 11).Value = 2400  # CUL!K97
$ws.Cells.Item(97, 12).Value = 2935.5  # CUL!L97
$ws.Cells.Item(97, 13).Value = -1904  # CUL!M97
$ws.Cells.Item(97, 14).Value = -3927.5  # CUL!N97

$ws.Cells.Item(113, 8).Value = 4701380.5  # CUL!H113
$ws.Cells.Item(113, 9).Value = 20833760  # CUL!I113
$ws.Cells.Item(113, 10).Value = 2857680  # CUL!J113
$ws.Cells.Item(113, 11).Value = 62501280  # CUL!K113
$ws.Cells.Item(113, 12).Value = 8573040  # CUL!L113
$ws.Cells.Item(113, 13).Value = -62499110  # CUL!M113

$ws.Cells.Item(126, 8).Value = 2708.3  # CUL!H126
$ws.Cells.Item(126, 9).Value = 837.5  # CUL!I126
$ws.Cells.Item(126, 10).Value = 2996.1155  # CUL!J126
$ws.Cells.Item(126, 11).Value = 2512.5  # CUL!K126
$ws.Cells.Item(126, 12).Value = 8988.3465  # CUL!L126
$ws.Cells.Item(126, 13).Value = 2427.5  # CUL!M126
$ws.Cells.Item(126, 14).Value = -18868.3465  # CUL!N126

$ws.Cells.Item(134, 8).Value = 41668640  # CUL!H134
$ws.Cells.Item(134, 9).Value = 41668640  # CUL!I134
$ws.Cells.Item(134, 10).Value = 0  # CUL!J134
$ws.Cells.Item(134, 11).Value = 125005920  # CUL!K134
$ws.Cells.Item(134, 12).Value = 0  # CUL!L134
$ws.Cells.Item(134, 13).Value = -125000850  # CUL!M134

$ws.Cells.Item(140, 8).Value = 6758914.5  # CUL!H140
$ws.Cells.Item(140, 9).Value = 15626659  # CUL!I140
$ws.Cells.Item(140, 10).Value = 2538.0476  # CUL!J140
$ws.Cells.Item(140, 11).Value = 46879977  # CUL!K140
$ws.Cells.Item(140, 12).Value = 7614.1428  # CUL!L140
$ws.Cells.Item(140, 13).Value = -46874797  # CUL!M140
$ws.Cells.Item(140, 14).Value = -17974.1428  # CUL!N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 14290400  # GSM!H80
$ws.Cells.Item(80, 9).Value = 5959.8  # GSM!I80
$ws.Cells.Item(80, 10).Value = 50001500  # GSM!J80
$ws.Cells.Item(80, 11).Value = 5959.8  # GSM!K80
$ws.Cells.Item(80, 12).Value = 50001500  # GSM!L80
$ws.Cells.Item(80, 13).Value = -4961.8  # GSM!M80
$ws.Cells.Item(80, 14).Value = -50003496  # GSM!N80

$ws.Cells.Item(83, 8).Value = 14290400  # GSM!H83
$ws.Cells.Item(83, 9).Value = 5959.8  # GSM!I83
$ws.Cells.Item(83, 10).Value = 50001500  # GSM!J83
$ws.Cells.Item(83, 11).Value = 29799  # GSM!K83
$ws.Cells.Item(83, 12).Value = 250007500  # GSM!L83
$ws.Cells.Item(83, 13).Value = -24807  # GSM!M83
$ws.Cells.Item(83, 14).Value = -250017484  # GSM!N83

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(129, 8).Value = 62000  # LTW!H129
$ws.Cells.Item(129, 9).Value = 0  # LTW!I129
$ws.Cells.Item(129, 10).Value = 62000  # LTW!J129
$ws.Cells.Item(129, 11).Value = 0  # LTW!K129
$ws.Cells.Item(129, 12).Value = 62000  # LTW!L129
$ws.Cells.Item(129, 14).Value = -72000  # LTW!N129

$ws.Cells.Item(133, 8).Value = 51884  # LTW!H133
$ws.Cells.Item(133, 9).Value = 0  # LTW!I133
$ws.Cells.Item(133, 10).Value = 51884  # LTW!J133
$ws.Cells.Item(133, 11).Value = 0  # LTW!K133
$ws.Cells.Item(133, 12).Value = 51884  # LTW!L133
$ws.Cells.Item(133, 14).Value = -56944  # LTW!N133

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 17394  # WVR!H132
$ws.Cells.Item(132, 9).Value = 20326.66  # WVR!I132
$ws.Cells.Item(132, 10).Value = 7679.5625  # WVR!J132
$ws.Cells.Item(132, 11).Value = 60979.98  # WVR!K132
$ws.Cells.Item(132, 12).Value = 23038.6875  # WVR!L132
$ws.Cells.Item(132, 13).Value = -58449.98  # WVR!M132
$ws.Cells.Item(132, 14).Value = -28098.6875  # WVR!N132
